$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Write the new data block for rows 2-5 (A:AH), replacing the old values
$arr = New-Object 'object[,]' 4,34
$arr[0,0] = 45043.50694444445
$arr[0,1] = 13.798
$arr[0,2] = 9.137
$arr[0,3] = 3.527
$arr[0,4] = 29.879
$arr[0,5] = 22.444
$arr[0,6] = 10.657
$arr[0,7] = 31.967
$arr[0,8] = 16.87
$arr[0,9] = 6.742
$arr[0,10] = 10.011
$arr[0,11] = 11.733
$arr[0,12] = 12.516
$arr[0,13] = 3.497
$arr[0,14] = 10.903
$arr[0,15] = 14.966
$arr[0,16] = 9.704000000000001
$arr[0,17] = 3.059
$arr[0,18] = 1.672
$arr[0,19] = 158.575
$arr[0,20] = 30.182
$arr[0,21] = 10.064
$arr[0,22] = 19.331
$arr[0,23] = 9.888999999999999
$arr[0,24] = 2.87
$arr[0,25] = 17.102
$arr[0,26] = 8.888999999999999
$arr[0,27] = 8.15
$arr[0,28] = 9.673999999999999
$arr[0,29] = 12.058
$arr[0,30] = 3.072
$arr[0,31] = 28.934
$arr[0,32] = 5.407
$arr[0,33] = 12.581
$arr[1,0] = 45043.51388888889
$arr[1,1] = 5.654
$arr[1,2] = 3.745
$arr[1,3] = 1.4
$arr[1,4] = 12.516
$arr[1,5] = 9.132
$arr[1,6] = 4.328
$arr[1,7] = 19.479
$arr[1,8] = 6.981
$arr[1,9] = 2.808
$arr[1,10] = 3.898
$arr[1,11] = 4.934
$arr[1,12] = 5.327
$arr[1,13] = 1.455
$arr[1,14] = 4.512
$arr[1,15] = 6.2
$arr[1,16] = 4.24
$arr[1,17] = 1.36
$arr[1,18] = 0.705
$arr[1,19] = 61.371
$arr[1,20] = 12.791
$arr[1,21] = 4.164
$arr[1,22] = 8.090999999999999
$arr[1,23] = 4.121
$arr[1,24] = 1.213
$arr[1,25] = 9.545999999999999
$arr[1,26] = 3.678
$arr[1,27] = 3.483
$arr[1,28] = 4.107
$arr[1,29] = 5.007
$arr[1,30] = 1.198
$arr[1,31] = 18.384
$arr[1,32] = 2.154
$arr[1,33] = 5.208
$arr[2,0] = 45043.52083333334
$arr[2,1] = 8.552
$arr[2,2] = 6.135
$arr[2,3] = 1.041
$arr[2,4] = 18.823
$arr[2,5] = 14.722
$arr[2,6] = 6.646
$arr[2,7] = 25.444
$arr[2,8] = 10.471
$arr[2,9] = 4.464
$arr[2,10] = 6.464
$arr[2,11] = 7.508
$arr[2,12] = 8.042
$arr[2,13] = 2.172
$arr[2,14] = 6.767
$arr[2,15] = 9.460000000000001
$arr[2,16] = 5.976
$arr[2,17] = 0.93
$arr[2,18] = 0.582
$arr[2,19] = 95.65900000000001
$arr[2,20] = 18.876
$arr[2,21] = 6.247
$arr[2,22] = 12.366
$arr[2,23] = 6.476
$arr[2,24] = 1.298
$arr[2,25] = 12.413
$arr[2,26] = 5.517
$arr[2,27] = 5.026
$arr[2,28] = 5.91
$arr[2,29] = 7.782
$arr[2,30] = 0.746
$arr[2,31] = 23.209
$arr[2,32] = 3.397
$arr[2,33] = 7.809
$arr[3,0] = 45043.52777777778
$arr[3,1] = 16.73
$arr[3,2] = 12.38
$arr[3,3] = 1.11
$arr[3,4] = 36.59
$arr[3,5] = 29.61
$arr[3,6] = 13.1
$arr[3,7] = 48.87
$arr[3,8] = 20.36
$arr[3,9] = 8.970000000000001
$arr[3,10] = 13.22
$arr[3,11] = 14.66
$arr[3,12] = 15.6
$arr[3,13] = 4.23
$arr[3,14] = 13.16
$arr[3,15] = 18.66
$arr[3,16] = 11.19
$arr[3,17] = 0.78
$arr[3,18] = 0.73
$arr[3,19] = 192.96
$arr[3,20] = 36.74
$arr[3,21] = 12.15
$arr[3,22] = 24.58
$arr[3,23] = 12.9
$arr[3,24] = 2.03
$arr[3,25] = 24.15
$arr[3,26] = 10.73
$arr[3,27] = 9.56
$arr[3,28] = 11.24
$arr[3,29] = 15.33
$arr[3,30] = 0.54
$arr[3,31] = 44.44
$arr[3,32] = 6.79
$arr[3,33] = 15.19
$ws.Range("A2:AH5").Value2 = $arr

# Step 2: Delete old row 6 (data no longer present in the new dataset)
$ws.Rows.Item(6).Delete()

# Step 3: Adjust column widths (raw OOXML width 7 -> 8) for columns G,K,L,M,O,P,V,AD,AH
$ws.Columns.Item(7).ColumnWidth = 7.17  # G: raw width 7 -> 8
$ws.Columns.Item(11).ColumnWidth = 7.17  # K: raw width 7 -> 8
$ws.Columns.Item(12).ColumnWidth = 7.17  # L: raw width 7 -> 8
$ws.Columns.Item(13).ColumnWidth = 7.17  # M: raw width 7 -> 8
$ws.Columns.Item(15).ColumnWidth = 7.17  # O: raw width 7 -> 8
$ws.Columns.Item(16).ColumnWidth = 7.17  # P: raw width 7 -> 8
$ws.Columns.Item(22).ColumnWidth = 7.17  # V: raw width 7 -> 8
$ws.Columns.Item(30).ColumnWidth = 7.17  # AD: raw width 7 -> 8
$ws.Columns.Item(34).ColumnWidth = 7.17  # AH: raw width 7 -> 8
